$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing last pair of columns: BH = 60 ("06-17_A"), BI = 61 ("06-17_0")
# New pair being appended:       BJ = 62 ("06-18_A"), BK = 63 ("06-18_0")
$colBH = 60
$colBI = 61
$colBJ = 62
$colBK = 63

$lastRow = 162

# --- Header row (row 1): style matches the existing header cells (bold, thin box border, centered/top) ---
$hdrBJ = $ws.Cells.Item(1, $colBJ)
$hdrBJ.Value = "06-18_A"
$hdrBJ.Font.Bold = $true
$hdrBJ.BorderAround(1, 2)
$hdrBJ.HorizontalAlignment = -4108
$hdrBJ.VerticalAlignment = -4160

$hdrBK = $ws.Cells.Item(1, $colBK)
$hdrBK.Value = "06-18_0"
$hdrBK.Font.Bold = $true
$hdrBK.BorderAround(1, 2)
$hdrBK.HorizontalAlignment = -4108
$hdrBK.VerticalAlignment = -4160

# --- Data rows (2..162) ---
for ($r = 2; $r -le $lastRow; $r++) {
    $bhCell = $ws.Cells.Item($r, $colBH)
    $biCell = $ws.Cells.Item($r, $colBI)
    $bhVal = $bhCell.Value()

    if (-not [string]::IsNullOrEmpty($bhVal)) {
        # Row has data: duplicate BH -> BJ (value + fill color); duplicate old BI text -> BK
        $biTextVal = $biCell.Value()

        $bjCell = $ws.Cells.Item($r, $colBJ)
        $bjCell.Value = $bhVal
        $bjCell.Interior.Color = $bhCell.Interior.Color

        $bkCell = $ws.Cells.Item($r, $colBK)
        $bkCell.Value = "'" + $biTextVal

        # Old BI cell: was stored as text; normalize to the equivalent numeric value
        $biCell.Value = [double]$biTextVal
    }
    else {
        # Row has no data: just carry BH's fill onto BJ; BK/BI stay empty
        $bjCell = $ws.Cells.Item($r, $colBJ)
        $bjCell.Interior.Color = $bhCell.Interior.Color
    }
}

# A couple of ID cells (column A) were stored as text in the source data;
# normalize them to numeric, matching the rest of the column.
$a161 = $ws.Cells.Item(161, 1)
$a161.Value = [double]$a161.Value()
$a162 = $ws.Cells.Item(162, 1)
$a162.Value = [double]$a162.Value()

Write-Host "Added 06-18_A/06-18_0 columns (BJ/BK) and normalized BI column."
